# "Added WTA calculation" — update the Threshold value in D5 and tidy up
# the header cell's border formatting, then leave the selection on D5
# (matching the author's final cursor position when they saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the WTA threshold value in D5 (new shared string "[89, 85, 82]")
$ws.Range("D5").Value = "[89, 85, 82]"

# Re-apply the border-only formatting to the header cell D1 (drops the
# redundant "apply fill" flag that a stray fill-paint had left behind).
$ws.Range("D1").Borders.LineStyle = 1

# Leave the selection on D5, matching the saved state of the workbook.
$ws.Range("D5").Select()
